$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 522.25
$ws.Range("I98").Value = 488.8889
$ws.Range("J98").Value = 702.4
$ws.Range("K98").Value = 488.8889
$ws.Range("L98").Value = 702.4
$ws.Range("M98").Value = 1009.1111
$ws.Range("N98").Value = -3698.4

$ws.Range("H100").Value = 22225544
$ws.Range("I100").Value = 50002000
$ws.Range("J100").Value = 4380
$ws.Range("K100").Value = 50002000
$ws.Range("L100").Value = 4380
$ws.Range("M100").Value = -50001459
$ws.Range("N100").Value = -5462

$ws.Range("H106").Value = 5054.909
$ws.Range("I106").Value = 5260.4
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 5260.4
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -4629.4
$ws.Range("N106").Value = -4262

$ws.Range("H107").Value = 66666876
$ws.Range("I107").Value = 83333544
$ws.Range("J107").Value = 186.66667
$ws.Range("K107").Value = 83333544
$ws.Range("L107").Value = 186.66667
$ws.Range("M107").Value = -83331624
$ws.Range("N107").Value = -4026.66667

$ws.Range("H122").Value = 522.25
$ws.Range("I122").Value = 488.8889
$ws.Range("J122").Value = 702.4
$ws.Range("K122").Value = 1466.6667
$ws.Range("L122").Value = 2107.2
$ws.Range("M122").Value = 983.3333
$ws.Range("N122").Value = -7007.2

$ws.Range("H128").Value = 6316.8423
$ws.Range("J128").Value = 6316.8423
$ws.Range("L128").Value = 6316.8423
$ws.Range("N128").Value = -16276.8423

$ws.Range("H132").Value = 6650.2607
$ws.Range("I132").Value = 9079.286
$ws.Range("J132").Value = 2871.7778
$ws.Range("K132").Value = 27237.858
$ws.Range("L132").Value = 8615.3334
$ws.Range("M132").Value = -24707.858
$ws.Range("N132").Value = -13675.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11667.116
$ws.Range("I32").Value = 2875.6458
$ws.Range("J32").Value = 31761.904
$ws.Range("K32").Value = 2875.6458
$ws.Range("L32").Value = 31761.904
$ws.Range("M32").Value = -2588.6458
$ws.Range("N32").Value = -32335.904

$ws.Range("H102").Value = 125002310
$ws.Range("I102").Value = 125002310
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 125002310
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -125000688
$ws.Range("N102").ClearContents()

$ws.Range("H110").Value = 6206.5835
$ws.Range("I110").Value = 9233.929
$ws.Range("J110").Value = 1968.3
$ws.Range("K110").Value = 9233.929
$ws.Range("L110").Value = 1968.3
$ws.Range("M110").Value = -7188.929
$ws.Range("N110").Value = -6058.3

$ws.Range("H121").Value = 32770.8
$ws.Range("J121").Value = 32770.8
$ws.Range("L121").Value = 32770.8
$ws.Range("N121").Value = -36264.8

$ws.Range("H122").Value = 1797
$ws.Range("I122").Value = 1446.625
$ws.Range("K122").Value = 4339.875
$ws.Range("M122").Value = -1889.875

$ws.Range("H124").Value = 37835.8
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 37835.8
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 37835.8
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -47655.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 27000
$ws.Range("J55").Value = 27000
$ws.Range("L55").Value = 27000
$ws.Range("N55").Value = -27546

$ws.Range("H105").Value = 2934.2104
$ws.Range("I105").Value = 2626.923
$ws.Range("J105").Value = 3600
$ws.Range("K105").Value = 2626.923
$ws.Range("L105").Value = 3600
$ws.Range("M105").Value = -879.9229999999998
$ws.Range("N105").Value = -7094

$ws.Range("H107").Value = 2329.3845
$ws.Range("I107").Value = 1966.8889
$ws.Range("J107").Value = 3145
$ws.Range("K107").Value = 1966.8889
$ws.Range("L107").Value = 3145
$ws.Range("M107").Value = -46.88889999999992
$ws.Range("N107").Value = -6985

$ws.Range("H123").Value = 24000
$ws.Range("J123").Value = 24000
$ws.Range("L123").Value = 24000
$ws.Range("N123").Value = -33800

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 44840
$ws.Range("J138").Value = 44840
$ws.Range("L138").Value = 44840
$ws.Range("N138").Value = -55120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1218.9166
$ws.Range("I58").Value = 807.5909
$ws.Range("J58").Value = 1865.2858
$ws.Range("K58").Value = 807.5909
$ws.Range("L58").Value = 1865.2858
$ws.Range("M58").Value = -604.5909
$ws.Range("N58").Value = -2271.2858

$ws.Range("H107").Value = 17242194
$ws.Range("I107").Value = 21739964
$ws.Range("J107").Value = 745
$ws.Range("K107").Value = 21739964
$ws.Range("L107").Value = 745
$ws.Range("M107").Value = -21738044
$ws.Range("N107").Value = -4585

$ws.Range("H122").Value = 3922376.2
$ws.Range("I122").Value = 9524391
$ws.Range("J122").Value = 965.8
$ws.Range("K122").Value = 28573173
$ws.Range("L122").Value = 2897.4
$ws.Range("M122").Value = -28570723
$ws.Range("N122").Value = -7797.4

$ws.Range("H134").Value = 4467.7646
$ws.Range("I134").Value = 5244.9165
$ws.Range("J134").Value = 2602.6
$ws.Range("K134").Value = 15734.7495
$ws.Range("L134").Value = 7807.799999999999
$ws.Range("M134").Value = -13199.7495
$ws.Range("N134").Value = -12877.8

$ws.Range("H136").Value = 1218.9166
$ws.Range("I136").Value = 807.5909
$ws.Range("J136").Value = 1865.2858
$ws.Range("K136").Value = 2422.7727
$ws.Range("L136").Value = 5595.857400000001
$ws.Range("M136").Value = 127.2273
$ws.Range("N136").Value = -10695.8574

$ws.Range("H138").Value = 21156
$ws.Range("J138").Value = 21156
$ws.Range("L138").Value = 21156
$ws.Range("N138").Value = -31436

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 8619
$ws.Range("I120").Value = 4980
$ws.Range("J120").Value = 9138.857
$ws.Range("K120").Value = 14940
$ws.Range("L120").Value = 27416.571
$ws.Range("M120").Value = -10102
$ws.Range("N120").Value = -37092.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3020.1936
$ws.Range("I122").Value = 2236.5
$ws.Range("J122").Value = 3665.5881
$ws.Range("K122").Value = 6709.5
$ws.Range("L122").Value = 10996.7643
$ws.Range("M122").Value = -4259.5
$ws.Range("N122").Value = -15896.7643

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 196.36842
$ws.Range("I55").Value = 45.4
$ws.Range("J55").Value = 762.5
$ws.Range("K55").Value = 45.4
$ws.Range("L55").Value = 762.5
$ws.Range("M55").Value = 127.6
$ws.Range("N55").Value = -1108.5

$ws.Range("H100").Value = 1466.6666
$ws.Range("I100").Value = 1466.6666
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1466.6666
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -925.6666
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 4702.65
$ws.Range("I122").Value = 3713.1667
$ws.Range("J122").Value = 6186.875
$ws.Range("K122").Value = 11139.5001
$ws.Range("L122").Value = 18560.625
$ws.Range("M122").Value = -8689.500100000001
$ws.Range("N122").Value = -23460.625

$ws.Range("H139").Value = 38827.668
$ws.Range("J139").Value = 38827.668
$ws.Range("L139").Value = 38827.668
$ws.Range("N139").Value = -49107.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 35716336
$ws.Range("I122").Value = 47620924
$ws.Range("J122").Value = 2572.1428
$ws.Range("K122").Value = 142862772
$ws.Range("L122").Value = 7716.428400000001
$ws.Range("M122").Value = -142860322
$ws.Range("N122").Value = -12616.4284
